# Remove open burning and natural emissions
#
# The "map" sheet maps inv_sector/scaling_sector to ceds_sector. The rows
# covering open-burning and natural emissions (3F Agricultural-residue
# burning, 11B Forest-fires, 11C Other-natural x2) should no longer be
# scaled together with a CEDS sector -- instead they get flagged with a
# new "Notes" column explaining they should be excluded from the CMIP6
# data product.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# Add the new "Notes" column header
$ws.Range("D1").Value() = "Notes"

$note = "Don't include in CMIP6 data product"

# 4F (row 35): 3F_Agricultural-residue-burning-on-fields
$ws.Range("B35").ClearContents()
$ws.Range("D35").Value() = $note

# 5A (row 36): 11B_Forest-fires
$ws.Range("B36").ClearContents()
$ws.Range("D36").Value() = $note

# 5C (row 37): 11C_Other-natural
$ws.Range("B37").ClearContents()
$ws.Range("D37").Value() = $note

# 5D (row 38): 11C_Other-natural
$ws.Range("B38").ClearContents()
$ws.Range("D38").Value() = $note

# Restore view state (best effort; selection is what the runtime persists)
$ws.Activate()
$ws.Range("B39").Select()
